# Anonymize "fedcore" -> "approach" and fix up the header-row border
# styling / stray empty cells on both sheets, per the commit:
#   "update of results and scripts. Anonimyzed fedcore"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

# Give the two empty header cells under the merged "original|fedcore|change"
# band their own light border (top+bottom on C1, top+bottom+right on D1)
# instead of the inherited bold/boxed header style.
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142

# Anonymize the "fedcore" column label.
$ws1.Range("C2").Value = "approach"

# D5 was a stray empty inline-string cell - remove it entirely.
$ws1.Range("D5").ClearContents()

# ---------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

$c1b = $ws2.Range("C1")
$c1b.ClearFormats()
$c1b.Borders.LineStyle = 1
$c1b.Borders.Item(7).LineStyle = -4142
$c1b.Borders.Item(10).LineStyle = -4142

$d1b = $ws2.Range("D1")
$d1b.ClearFormats()
$d1b.Borders.LineStyle = 1
$d1b.Borders.Item(7).LineStyle = -4142

$xlPasteFormats = -4122

# F1/G1 need exactly the same two styles as C1/D1 - copy them over instead
# of re-deriving the border sequence, so no throwaway intermediate styles
# get left behind in the shared style table.
$c1b.Copy()
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$d1b.Copy()
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Anonymize both "fedcore" column labels.
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 was a stray empty inline-string cell - remove it entirely.
$ws2.Range("G5").ClearContents()
